# Update automatico via Actualizar 05-14-2020 16-03-32
#
# The "Fecha" column (A) is refreshed with a new run of dates: most rows keep
# a real date serial, but a block of rows (the ones for which no serial date
# was available in the refreshed export) switch to a plain text date such as
# "13/4/2020". Those text cells, plus the very last row, also pick up a new
# cell style that right-aligns the value. A handful of rows in the lower half
# of the table gain values in column C ("Pruebas Positivas ") that were
# missing before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column A ("Fecha") refresh, row by row (top to bottom, so that the
#    new text dates land in sharedStrings.xml in the same order Excel
#    would naturally append them).
# ---------------------------------------------------------------------
$ws.Range("A5").Value  = 43834
$ws.Range("A6").Value  = 43865
$ws.Range("A7").Value  = 43894
$ws.Range("A8").Value  = 43955
$ws.Range("A9").Value  = 43986
$ws.Range("A10").Value = 44016
$ws.Range("A11").Value = 44047
$ws.Range("A12").Value = 44078
$ws.Range("A13").Value = 44108
$ws.Range("A14").Value = 44139
$ws.Range("A15").Value = 44169

$ws.Range("A16").Value = "13/4/2020"
$ws.Range("A17").Value = "14/4/2020"
$ws.Range("A18").Value = "15/4/2020"
$ws.Range("A19").Value = "16/4/2020"
$ws.Range("A20").Value = "17/4/2020"
$ws.Range("A21").Value = "19/4/2020"
$ws.Range("A22").Value = "20/4/2020"
$ws.Range("A23").Value = "21/4/2020"
$ws.Range("A24").Value = "22/4/2020"
$ws.Range("A25").Value = "23/4/2020"
$ws.Range("A26").Value = "24/4/2020"
$ws.Range("A27").Value = "25/4/2020"
$ws.Range("A28").Value = "26/4/2020"
$ws.Range("A29").Value = "27/4/2020"
$ws.Range("A30").Value = "28/4/2020"
$ws.Range("A31").Value = "29/4/2020"
$ws.Range("A32").Value = "30/4/2020"

$ws.Range("A33").Value = 43835
$ws.Range("A34").Value = 43866
$ws.Range("A35").Value = 43895
$ws.Range("A36").Value = 43926
$ws.Range("A37").Value = 43956
$ws.Range("A38").Value = 43987
$ws.Range("A39").Value = 44017
$ws.Range("A40").Value = 44048
$ws.Range("A41").Value = 44079
$ws.Range("A42").Value = 44109
$ws.Range("A43").Value = 44140
$ws.Range("A44").Value = 44170

$ws.Range("A45").Value = "13/5/2020"

# ---------------------------------------------------------------------
# 2) Re-style rows 16-45 of column A: keep the date number format but
#    right-align the cell (new cellXf, applied as one batch so every row
#    shares the same style entry).
# ---------------------------------------------------------------------
$ws.Range("A16:A45").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 3) Fill in the previously-missing "Pruebas Positivas " (column C)
#    figures for the most recent rows.
# ---------------------------------------------------------------------
$ws.Range("C37").Value = 33
$ws.Range("C38").Value = 35
$ws.Range("C39").Value = 34
$ws.Range("C43").Value = 62
$ws.Range("C44").Value = 85
$ws.Range("C45").Value = 143

# ---------------------------------------------------------------------
# 4) Reflect the author's scroll position / active selection when they
#    saved the workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C48").Select() | Out-Null
